$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row (row 43) with the 2020-03-30 update
$ws.Cells.Item(43, 1).Value = 43920
$ws.Cells.Item(43, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(43, 2).Value = 6534
$ws.Cells.Item(43, 3).Value = 1414
$ws.Cells.Item(43, 4).Value = 122
$ws.Cells.Item(43, 5).Value = 75
$ws.Cells.Item(43, 6).Value = 5249
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 1107
$ws.Cells.Item(43, 9).Value = "(dari situs)"

# Move the active selection to G39 (matches the updated view state)
$ws.Range("G39").Select() | Out-Null
